# Fixes issues resulting from extra ticker rows being mixed into the
# MICROSOFT_CORP sheet: every row is re-pointed at the MSFT ticker and
# its price/volume columns (D:H) are corrected to the real MSFT series.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> [open_price, close_price, high_price, low_price, shares_outstanding]
$data = @{
  2 = @(34.91915977733629, 41.83418655395508, 42.60825793974698, 34.50632279855949, 7433087554)
  3 = @(38.4875771892215, 40.42667388916016, 41.03264181428221, 37.50071674369165, 7433087554)
  4 = @(38.99406957329596, 45.86922454833984, 47.37670438100428, 38.12269371690945, 7433087554)
  5 = @(47.65218106813233, 48.32766342163086, 48.59083753880731, 43.07293871278824, 7433087554)
  6 = @(48.63932551282061, 44.06254577636719, 50.1590294259118, 43.60310032891011, 7433087554)
  7 = @(45.49178630738633, 50.42977523803711, 50.97250977419564, 44.83338619941309, 7433087554)
  8 = @(51.39763721796317, 53.64477157592773, 54.94291847328731, 50.42178923748696, 7433087554)
  9 = @(56.59396941249549, 58.27042770385742, 59.40609455077209, 55.83685932730842, 7433087554)
  10 = @(59.67555393237211, 62.07853698730469, 62.69515142874145, 58.80504064767749, 7433087554)
  11 = @(63.22779167347375, 66.3011703491211, 67.86978243607513, 62.03308904936114, 7433087554)
  12 = @(68.49727068807812, 76.26292419433594, 79.03178418727094, 67.58042861168309, 7433087554)
  13 = @(79.36418172743693, 87.546630859375, 87.95206252468533, 78.78367298364888, 7433087554)
  14 = @(83.75485902486797, 86.5784683227539, 90.6333670354583, 81.01456608037077, 7433087554)
  15 = @(91.21216563133493, 98.63187408447266, 103.3458971130981, 91.11918828652058, 7433087554)
  16 = @(107.1037212317319, 99.69279479980467, 108.4384345567588, 93.43924704034512, 7433087554)
  17 = @(93.31797174683828, 97.89247131347656, 101.1452434475639, 91.115080772955, 7433087554)
  18 = @(111.9796871890666, 122.947021484375, 123.6718901779615, 111.1794974817505, 7433087554)
  19 = @(129.0997725573902, 128.7596130371094, 133.8714350926367, 127.2477897430132, 7433087554)
  20 = @(132.402210907896, 135.9194030761719, 138.0998805297197, 126.2968798242927, 7433087554)
  21 = @(151.040760577142, 161.9326629638672, 165.5664758880142, 148.8814006170601, 7433087554)
  22 = @(145.9400560609984, 170.9406433105469, 172.0757204095209, 143.4218752837023, 7433087554)
  23 = @(194.3060821468002, 196.0947570800781, 206.970321600913, 188.9209084019272, 7433087554)
  24 = @(204.6996530148239, 194.1333923339844, 215.9370885925192, 191.4007339168742, 7433087554)
  25 = @(213.9260070850278, 222.9914093017578, 233.2584663338894, 203.7454667352439, 7433087554)
  26 = @(229.777716845021, 242.9879760742188, 253.596668544704, 229.3730277023811, 7433087554)
  27 = @(260.3824460960315, 275.1588134765625, 280.2194623726587, 260.3728083928387, 7433087554)
  28 = @(272.9859115038698, 320.8832702636719, 321.2509718839253, 271.1764604532231, 7433087554)
  29 = @(325.0862853620566, 301.462158203125, 327.6551735575567, 267.6012031386585, 7433087554)
  30 = @(300.5215017665168, 269.5824584960937, 306.0973187900976, 262.2775536012272, 7433087554)
  31 = @(249.6368757403655, 273.3454895019531, 274.5723114658439, 239.4621097597596, 7433087554)
  32 = @(229.6960056147813, 226.495620727539, 244.946615735801, 213.8111637098241, 7433087554)
  33 = @(237.8482622712028, 242.4764556884766, 244.4529840016496, 214.629000279338, 7433087554)
  34 = @(281.0554951876676, 301.3999633789062, 303.038096464785, 270.1181533753086, 7433087554)
  35 = @(333.4482780270713, 330.233642578125, 360.5712377352722, 321.4646249300586, 7433087554)
  36 = @(311.5844666663645, 333.0903625488281, 341.0602838616102, 306.5897293000804, 7433087554)
  37 = @(369.0571049559907, 392.4723815917969, 409.9845009870256, 361.7916713905026, 7433087554)
  38 = @(419.2775635448957, 385.0390930175781, 424.6378116014886, 383.753432683152, 7433087554)
  39 = @(444.5155790150141, 414.4855651855469, 464.0236982246781, 408.4022679352166, 7433087554)
  40 = @(425.257410073319, 403.3220825195313, 435.2325102211543, 403.2724369212272, 7433087554)
  41 = @(423.2045320573043, 412.791748046875, 445.9296658158772, 408.475469230866, 7433087554)
  42 = @(373.3496987636053, 393.88818359375, 395.2833185411371, 343.5933482559051, 7433087554)
  43 = @(495.6551665959438, 532.6243896484375, 554.5383762548963, 487.8979301273739, 7433087554)
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  $ws.Cells.Item($r, 4).Value = $row[0]   # D: open_price
  $ws.Cells.Item($r, 5).Value = $row[1]   # E: close_price
  $ws.Cells.Item($r, 6).Value = $row[2]   # F: high_price
  $ws.Cells.Item($r, 7).Value = $row[3]   # G: low_price
  $ws.Cells.Item($r, 8).Value = $row[4]   # H: shares_outstanding
  $ws.Cells.Item($r, 9).Value = "MSFT"    # I: fixed_ticker
}
